# Rename the inline-picture "name" metadata for the two repeated logos.
#
#   Pearson logo (footers): wp:docPr/@name  "image2.png" -> "image1.png"
#   BTec logo   (headers):  wp:docPr/@name  "image1.jpg" -> "image2.jpg"
#
# The underlying image bytes / relationship ids are untouched - only the
# docPr "name" attribute of each <wp:inline> picture changes.

$d = $word.ActiveDocument
$sec = $d.Sections(1)

# --- Headers: BTec_Logo-Orange picture -> name="image2.jpg" -----------
for ($i = 1; $i -le 3; $i++) {
    $hdr = $sec.Headers($i)
    if ($hdr.Exists) {
        $cnt = $hdr.Range.InlineShapes.Count
        for ($j = 1; $j -le $cnt; $j++) {
            $shp = $hdr.Range.InlineShapes($j)
            if ($shp.AlternativeText -eq "BTec_Logo-Orange") {
                $shp.Name = "image2.jpg"
            }
        }
    }
}

# --- Footers: Pearson logo picture -> name="image1.png" ---------------
for ($i = 1; $i -le 3; $i++) {
    $ftr = $sec.Footers($i)
    if ($ftr.Exists) {
        $cnt = $ftr.Range.InlineShapes.Count
        for ($j = 1; $j -le $cnt; $j++) {
            $shp = $ftr.Range.InlineShapes($j)
            if ($shp.AlternativeText -eq "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png") {
                $shp.Name = "image1.png"
            }
        }
    }
}

Write-Output "done"
